{"js": "const pairs = [\n  [\"14\u00d715=210\", \"31\u00d791=2821\"],\n  [\"55\u00d772=3960\", \"42\u00d713=546\"],\n  [\"55\u00d760=3300\", \"59\u00d753=3127\"],\n  [\"26\u00d776=1976\", \"88\u00d787=7656\"],\n  [\"51\u00d737=1887\", \"97\u00d768=6596\"],\n  [\"27\u00d780=2160\", \"54\u00d740=2160\"],\n  [\"74\u00d795=7030\", \"27\u00d729=783\"],\n  [\"40\u00d712=480\", \"53\u00d764=3392\"],\n  [\"95\u00d734=3230\", \"98\u00d780=7840\"],\n  [\"84\u00d791=7644\", \"27\u00d759=1593\"],\n  [\"15\u00d756=840\", \"44\u00d792=4048\"],\n  [\"62\u00d726=1612\", \"79\u00d716=1264\"],\n  [\"24\u00d733=792\", \"13\u00d723=299\"],\n  [\"67\u00d724=1608\", \"79\u00d787=6873\"],\n  [\"43\u00d716=688\", \"96\u00d753=5088\"],\n  [\"22\u00d783=1826\", \"95\u00d783=7885\"],\n  [\"90\u00d777=6930\", \"20\u00d717=340\"],\n  [\"77\u00d743=3311\", \"54\u00d739=2106\"],\n  [\"60\u00d755=3300\", \"20\u00d718=360\"],\n  [\"65\u00d774=4810\", \"20\u00d744=880\"],\n  [\"21\u00d772=1512\", \"35\u00d746=1610\"],\n  [\"86\u00d790=7740\", \"14\u00d727=378\"],\n  [\"17\u00d750=850\", \"90\u00d713=1170\"],\n  [\"58\u00d734=1972\", \"89\u00d759=5251\"],\n  [\"49\u00d719=931\", \"55\u00d754=2970\"],\n];\n\nconst body = context.document.body;\nfor (const [before, after] of pairs) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(after, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @('14\u00d715=210', '31\u00d791=2821'),\n    @('55\u00d772=3960', '42\u00d713=546'),\n    @('55\u00d760=3300', '59\u00d753=3127'),\n    @('26\u00d776=1976', '88\u00d787=7656'),\n    @('51\u00d737=1887', '97\u00d768=6596'),\n    @('27\u00d780=2160', '54\u00d740=2160'),\n    @('74\u00d795=7030', '27\u00d729=783'),\n    @('40\u00d712=480', '53\u00d764=3392'),\n    @('95\u00d734=3230', '98\u00d780=7840'),\n    @('84\u00d791=7644', '27\u00d759=1593'),\n    @('15\u00d756=840', '44\u00d792=4048'),\n    @('62\u00d726=1612', '79\u00d716=1264'),\n    @('24\u00d733=792', '13\u00d723=299'),\n    @('67\u00d724=1608', '79\u00d787=6873'),\n    @('43\u00d716=688', '96\u00d753=5088'),\n    @('22\u00d783=1826', '95\u00d783=7885'),\n    @('90\u00d777=6930', '20\u00d717=340'),\n    @('77\u00d743=3311', '54\u00d739=2106'),\n    @('60\u00d755=3300', '20\u00d718=360'),\n    @('65\u00d774=4810', '20\u00d744=880'),\n    @('21\u00d772=1512', '35\u00d746=1610'),\n    @('86\u00d790=7740', '14\u00d727=378'),\n    @('17\u00d750=850', '90\u00d713=1170'),\n    @('58\u00d734=1972', '89\u00d759=5251'),\n    @('49\u00d719=931', '55\u00d754=2970'),\n)\n\nforeach ($pair in $replacements) {\n    $find = $pair[0]\n    $replace = $pair[1]\n    $range = $d.Content\n    $null = $range.Find.Execute(\n        $find,    # FindText\n        $false,   # MatchCase\n        $false,   # MatchWholeWord\n        $false,   # MatchWildcards\n        $false,   # MatchSoundsLike\n        $false,   # MatchAllWordForms\n        $true,    # Forward\n        1,        # Wrap (wdFindContinue)\n        $false,   # Format\n        $replace, # ReplaceWith\n        2         # Replace (wdReplaceAll)\n    )\n}\n"}
